# "Artefakty z Zad Domowe 1"
# Fill in the missing "Zadanie domowe 1" (Homework 1) score for the
# student in row 7 of Tabela1 on the "Sheet2" worksheet. Dependent
# table formulas (Suma pkt / % / Propozycja oceny) and the COUNTIF
# summary table below the table recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("M7").Value = 2.5001000000000002

# Recalculate so every dependent formula cell picks up the new value.
$excel.CalculateFullRebuild()
$excel.Calculate()

# Match the author's resulting view state: the sheet had scrolled so
# column C was left-most with M4 selected; afterwards the view is back
# to the default left edge (column A) with M7 selected.
$ws.Activate()
$ws.Range("M7").Select()
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
